$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor values change
$ws.Range("B3").Value = 0.0368061972743692
$ws.Range("C3").Value = 0.03802866507999527
$ws.Range("D3").Value = 0.0384258460334149

# Row 4 - label change + values change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03495736580196739
$ws.Range("C4").Value = 0.03503668499210437
$ws.Range("D4").Value = 0.03501326136718511

# Row 5 - label change + values change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.03614895606264568
$ws.Range("C5").Value = 0.03554943969873335
$ws.Range("D5").Value = 0.03268968492600135
